# Atualizado por script em 26-11-2023 20:30
# Adds 3 new match rows (36, 37, 38) to the India ISL 2023-2024 sheet,
# mirroring the formatting of the last existing data row (35).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting/styles of the last existing row down into the new rows
$srcRow = $ws.Range("A35:V35")
$srcRow.Copy($ws.Range("A36:V36"))
$srcRow.Copy($ws.Range("A37:V37"))
$srcRow.Copy($ws.Range("A38:V38"))

# ---- Row 36 ----
$ws.Cells.Item(36, 1).Value = 35
$ws.Cells.Item(36, 2).Value = "india"
$ws.Cells.Item(36, 3).Value = "isl"
$ws.Cells.Item(36, 4).Value = "2023-2024"
$ws.Cells.Item(36, 5).Value = 45255.54166666666
$ws.Cells.Item(36, 6).Value = "Chennaiyin"
$ws.Cells.Item(36, 7).Value = 1
$ws.Cells.Item(36, 8).Value = "East Bengal"
$ws.Cells.Item(36, 9).Value = 1
$ws.Cells.Item(36, 10).Value = 1.86
$ws.Cells.Item(36, 11).Value = "18/11/2023 13:13"
$ws.Cells.Item(36, 12).Value = 2.51
$ws.Cells.Item(36, 13).Value = "25/11/2023 12:59"
$ws.Cells.Item(36, 14).Value = 3.79
$ws.Cells.Item(36, 15).Value = "18/11/2023 13:13"
$ws.Cells.Item(36, 16).Value = 3.43
$ws.Cells.Item(36, 17).Value = "25/11/2023 12:59"
$ws.Cells.Item(36, 18).Value = 3.9
$ws.Cells.Item(36, 19).Value = "18/11/2023 13:13"
$ws.Cells.Item(36, 20).Value = 2.56
$ws.Cells.Item(36, 21).Value = "25/11/2023 12:59"
$ws.Cells.Item(36, 22).Value = "https://www.betexplorer.com/football/india/isl/chennaiyin-fc-east-bengal/G0ztV78q/"

# ---- Row 37 ----
$ws.Cells.Item(37, 1).Value = 36
$ws.Cells.Item(37, 2).Value = "india"
$ws.Cells.Item(37, 3).Value = "isl"
$ws.Cells.Item(37, 4).Value = "2023-2024"
$ws.Cells.Item(37, 5).Value = 45255.64583333334
$ws.Cells.Item(37, 6).Value = "Kerala Blasters"
$ws.Cells.Item(37, 7).Value = 1
$ws.Cells.Item(37, 8).Value = "Hyderabad"
$ws.Cells.Item(37, 9).Value = 0
$ws.Cells.Item(37, 10).Value = 2.52
$ws.Cells.Item(37, 11).Value = "18/11/2023 18:13"
$ws.Cells.Item(37, 12).Value = 1.93
$ws.Cells.Item(37, 13).Value = "25/11/2023 15:25"
$ws.Cells.Item(37, 14).Value = 3.15
$ws.Cells.Item(37, 15).Value = "18/11/2023 18:13"
$ws.Cells.Item(37, 16).Value = 3.38
$ws.Cells.Item(37, 17).Value = "25/11/2023 15:25"
$ws.Cells.Item(37, 18).Value = 2.94
$ws.Cells.Item(37, 19).Value = "18/11/2023 18:13"
$ws.Cells.Item(37, 20).Value = 4.26
$ws.Cells.Item(37, 21).Value = "25/11/2023 15:25"
$ws.Cells.Item(37, 22).Value = "https://www.betexplorer.com/football/india/isl/kerala-blasters-hyderabad/SzvpURNk/"

# ---- Row 38 ----
$ws.Cells.Item(38, 1).Value = 37
$ws.Cells.Item(38, 2).Value = "india"
$ws.Cells.Item(38, 3).Value = "isl"
$ws.Cells.Item(38, 4).Value = "2023-2024"
$ws.Cells.Item(38, 5).Value = 45256.64583333334
$ws.Cells.Item(38, 6).Value = "North East Utd"
$ws.Cells.Item(38, 7).Value = 1
$ws.Cells.Item(38, 8).Value = "Bengaluru FC"
$ws.Cells.Item(38, 9).Value = 1
$ws.Cells.Item(38, 10).Value = 3.36
$ws.Cells.Item(38, 11).Value = "19/11/2023 15:42"
$ws.Cells.Item(38, 12).Value = 2.87
$ws.Cells.Item(38, 13).Value = "26/11/2023 15:26"
$ws.Cells.Item(38, 14).Value = 3.63
$ws.Cells.Item(38, 15).Value = "19/11/2023 15:42"
$ws.Cells.Item(38, 16).Value = 3.63
$ws.Cells.Item(38, 17).Value = "26/11/2023 15:26"
$ws.Cells.Item(38, 18).Value = 2.12
$ws.Cells.Item(38, 19).Value = "19/11/2023 15:42"
$ws.Cells.Item(38, 20).Value = 2.38
$ws.Cells.Item(38, 21).Value = "26/11/2023 15:26"
$ws.Cells.Item(38, 22).Value = "https://www.betexplorer.com/football/india/isl/north-east-united-bengaluru-fc/AaTlTowd/"

$wb.Save()
